# Update crypto price/volume values (data refresh) on sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format first so Excel does not
# reinterpret numeric-looking strings (e.g. "10.60", "1.00") as numbers,
# which would silently drop significant trailing/format digits.
$updateRange = $ws.Range("D2:E51")
$updateRange.NumberFormat = "@"

$ws.Range("D2").Value = "57.290.82"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "2.349.32"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "520.26"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "136.13"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "2.361.48"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  +5.13%  "
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "2.767.39"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "57.274.00"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.355.90"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "328.72"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "61.24"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "8.28"
$ws.Range("E27").Value = "  +8.12%  "
$ws.Range("D28").Value = "1.32"
$ws.Range("E28").Value = "  +10.04%  "
$ws.Range("D29").Value = "170.35"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "18.57"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "1.30"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "0.922"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +3.96%  "
$ws.Range("D40").Value = "38.52"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("D41").Value = "150.81"
$ws.Range("E41").Value = "  +7.28%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "3.65"
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").Value = "281.95"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").Value = "0.0938"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "0.562"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "18.21"
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("E51").Value = "  +3.80%  "

# Restore the original (default / General) cell formatting so only the
# cell text content changes, matching the source data refresh.
$updateRange.ClearFormats()

